$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A343").Copy($ws.Range("A344"))
$ws.Range("A344").Value = 44418
$ws.Range("B344").Value = 0
$ws.Range("C344").Value = 3
$ws.Range("D344").Value = 71.47962830593281

$ws.Range("A343").Copy($ws.Range("A345"))
$ws.Range("A345").Value = 44419
$ws.Range("B345").Value = 0
$ws.Range("C345").Value = 3
$ws.Range("D345").Value = 71.47962830593281

$ws.Range("A343").Copy($ws.Range("A346"))
$ws.Range("A346").Value = 44420
$ws.Range("B346").Value = 2
$ws.Range("C346").Value = 4
$ws.Range("D346").Value = 95.30617107457708

$ws.Range("A343").Copy($ws.Range("A347"))
$ws.Range("A347").Value = 44421
$ws.Range("B347").Value = 1
$ws.Range("C347").Value = 4
$ws.Range("D347").Value = 95.30617107457708

$ws.Range("A343").Copy($ws.Range("A348"))
$ws.Range("A348").Value = 44422
$ws.Range("B348").Value = 0
$ws.Range("C348").Value = 4
$ws.Range("D348").Value = 95.30617107457708

$ws.Range("A343").Copy($ws.Range("A349"))
$ws.Range("A349").Value = 44423
$ws.Range("B349").Value = 0
$ws.Range("C349").Value = 3
$ws.Range("D349").Value = 71.47962830593281

$ws.Range("A343").Copy($ws.Range("A350"))
$ws.Range("A350").Value = 44424
$ws.Range("B350").Value = 1
$ws.Range("C350").Value = 4
$ws.Range("D350").Value = 95.30617107457708

$ws.Range("A343").Copy($ws.Range("A351"))
$ws.Range("A351").Value = 44425
$ws.Range("B351").Value = 1
$ws.Range("C351").Value = 5
$ws.Range("D351").Value = 119.1327138432213

$ws.Range("A343").Copy($ws.Range("A352"))
$ws.Range("A352").Value = 44426
$ws.Range("B352").Value = 0
$ws.Range("C352").Value = 5
$ws.Range("D352").Value = 119.1327138432213

$ws.Range("A343").Copy($ws.Range("A353"))
$ws.Range("A353").Value = 44427
$ws.Range("B353").Value = 0
$ws.Range("C353").Value = 3
$ws.Range("D353").Value = 71.47962830593281

$ws.Range("A343").Copy($ws.Range("A354"))
$ws.Range("A354").Value = 44428
$ws.Range("B354").Value = 1
$ws.Range("C354").Value = 3
$ws.Range("D354").Value = 71.47962830593281

$ws.Range("A343").Copy($ws.Range("A355"))
$ws.Range("A355").Value = 44429
$ws.Range("B355").Value = 0
$ws.Range("C355").Value = 3
$ws.Range("D355").Value = 71.47962830593281

$ws.Range("A343").Copy($ws.Range("A356"))
$ws.Range("A356").Value = 44430
$ws.Range("B356").Value = 0
$ws.Range("C356").Value = 3
$ws.Range("D356").Value = 71.47962830593281

$ws.Range("A343").Copy($ws.Range("A357"))
$ws.Range("A357").Value = 44431
$ws.Range("B357").Value = 1
$ws.Range("C357").Value = 3
$ws.Range("D357").Value = 71.47962830593281

